{"js": "const replacements = [\n  [\"26\u00d737=962\", \"17\u00d782=1394\"],\n  [\"24\u00d732=768\", \"18\u00d796=1728\"],\n  [\"71\u00d796=6816\", \"43\u00d723=989\"],\n  [\"14\u00d735=490\", \"36\u00d768=2448\"],\n  [\"27\u00d799=2673\", \"17\u00d780=1360\"],\n  [\"93\u00d746=4278\", \"53\u00d728=1484\"],\n  [\"66\u00d796=6336\", \"20\u00d741=820\"],\n  [\"46\u00d743=1978\", \"76\u00d733=2508\"],\n  [\"63\u00d752=3276\", \"61\u00d726=1586\"],\n  [\"48\u00d773=3504\", \"58\u00d793=5394\"],\n  [\"94\u00d767=6298\", \"85\u00d773=6205\"],\n  [\"60\u00d744=2640\", \"72\u00d716=1152\"],\n  [\"13\u00d787=1131\", \"20\u00d721=420\"],\n  [\"14\u00d736=504\", \"42\u00d759=2478\"],\n  [\"75\u00d722=1650\", \"50\u00d732=1600\"],\n  [\"16\u00d741=656\", \"86\u00d773=6278\"],\n  [\"60\u00d726=1560\", \"47\u00d724=1128\"],\n  [\"30\u00d793=2790\", \"96\u00d719=1824\"],\n  [\"24\u00d789=2136\", \"68\u00d794=6392\"],\n  [\"21\u00d750=1050\", \"60\u00d742=2520\"],\n  [\"85\u00d715=1275\", \"52\u00d771=3692\"],\n  [\"38\u00d738=1444\", \"66\u00d732=2112\"],\n  [\"38\u00d779=3002\", \"83\u00d778=6474\"],\n  [\"84\u00d716=1344\", \"51\u00d742=2142\"],\n  [\"58\u00d790=5220\", \"61\u00d777=4697\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"26\u00d737=962\", \"17\u00d782=1394\"),\n  @(\"24\u00d732=768\", \"18\u00d796=1728\"),\n  @(\"71\u00d796=6816\", \"43\u00d723=989\"),\n  @(\"14\u00d735=490\", \"36\u00d768=2448\"),\n  @(\"27\u00d799=2673\", \"17\u00d780=1360\"),\n  @(\"93\u00d746=4278\", \"53\u00d728=1484\"),\n  @(\"66\u00d796=6336\", \"20\u00d741=820\"),\n  @(\"46\u00d743=1978\", \"76\u00d733=2508\"),\n  @(\"63\u00d752=3276\", \"61\u00d726=1586\"),\n  @(\"48\u00d773=3504\", \"58\u00d793=5394\"),\n  @(\"94\u00d767=6298\", \"85\u00d773=6205\"),\n  @(\"60\u00d744=2640\", \"72\u00d716=1152\"),\n  @(\"13\u00d787=1131\", \"20\u00d721=420\"),\n  @(\"14\u00d736=504\", \"42\u00d759=2478\"),\n  @(\"75\u00d722=1650\", \"50\u00d732=1600\"),\n  @(\"16\u00d741=656\", \"86\u00d773=6278\"),\n  @(\"60\u00d726=1560\", \"47\u00d724=1128\"),\n  @(\"30\u00d793=2790\", \"96\u00d719=1824\"),\n  @(\"24\u00d789=2136\", \"68\u00d794=6392\"),\n  @(\"21\u00d750=1050\", \"60\u00d742=2520\"),\n  @(\"85\u00d715=1275\", \"52\u00d771=3692\"),\n  @(\"38\u00d738=1444\", \"66\u00d732=2112\"),\n  @(\"38\u00d779=3002\", \"83\u00d778=6474\"),\n  @(\"84\u00d716=1344\", \"51\u00d742=2142\"),\n  @(\"58\u00d790=5220\", \"61\u00d777=4697\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}"}
